$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.006876353814593728
$ws.Range("C2").Value = 87981.0709163148
$ws.Range("D2").Value = 10137753.70137369
$ws.Range("E2").Value = 5548678842208.939
$ws.Range("G2").Value = 5548689067943.719

# Row 3
$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 6.048734245549538
